# Swap the full data (all columns except the row-number column A) between
# pairs of rows: (161,162), (180,181) and (184,185).
#
# The underlying change is that two rows describing two different matches
# had been written to the wrong row, so this edit swaps their content back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $rowA, $rowB, $firstCol, $lastCol) {
    $rangeA = $ws.Range($ws.Cells.Item($rowA, $firstCol), $ws.Cells.Item($rowA, $lastCol))
    $rangeB = $ws.Range($ws.Cells.Item($rowB, $firstCol), $ws.Cells.Item($rowB, $lastCol))

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Column B is the first column (2), column AD is the last column (30).
# Column A (the row index) is intentionally left untouched.
$firstCol = 2   # B
$lastCol = 30   # AD

Swap-RowData $ws 161 162 $firstCol $lastCol
Swap-RowData $ws 180 181 $firstCol $lastCol
Swap-RowData $ws 184 185 $firstCol $lastCol
